$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1548.2417
$ws.Range("I15").Value = 1548.2417
$ws.Range("K15").Value = 4644.7251
$ws.Range("M15").Value = -4475.7251
$ws.Range("H132").Value = 2092.1562
$ws.Range("I132").Value = 2137.3872
$ws.Range("J132").Value = 690
$ws.Range("K132").Value = 6412.1616
$ws.Range("L132").Value = 2070
$ws.Range("M132").Value = -3882.1616
$ws.Range("N132").Value = -7130
$ws.Range("H134").Value = 53492.145
$ws.Range("J134").Value = 53492.145
$ws.Range("L134").Value = 53492.145
$ws.Range("N134").Value = -63632.145
$ws.Range("H136").Value = 74331.664
$ws.Range("J136").Value = 74331.664
$ws.Range("L136").Value = 74331.664
$ws.Range("N136").Value = -84531.664
$ws.Range("H137").Value = 2095.889
$ws.Range("I137").Value = 1665.7241
$ws.Range("J137").Value = 2875.5625
$ws.Range("K137").Value = 4997.1723
$ws.Range("L137").Value = 8626.6875
$ws.Range("M137").Value = -2447.1723
$ws.Range("N137").Value = -13726.6875
$ws.Range("H138").Value = 3409.8057
$ws.Range("I138").Value = 1625
$ws.Range("J138").Value = 8764.223
$ws.Range("K138").Value = 4875
$ws.Range("L138").Value = 26292.669
$ws.Range("M138").Value = 265
$ws.Range("N138").Value = -36572.669
$ws.Range("H139").Value = 67423.75
$ws.Range("J139").Value = 67423.75
$ws.Range("L139").Value = 67423.75
$ws.Range("N139").Value = -77703.75
$ws.Range("H140").Value = 124881.25
$ws.Range("J140").Value = 124881.25
$ws.Range("L140").Value = 124881.25
$ws.Range("N140").Value = -135241.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17972.973
$ws.Range("I32").Value = 19195.936
$ws.Range("K32").Value = 19195.936
$ws.Range("M32").Value = -18908.936
$ws.Range("H126").Value = 6600
$ws.Range("I126").Value = 6600
$ws.Range("K126").Value = 19800
$ws.Range("M126").Value = -17330

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1676.75
$ws.Range("I99").Value = 1501.1111
$ws.Range("J99").Value = 2203.6667
$ws.Range("K99").Value = 1501.1111
$ws.Range("L99").Value = 2203.6667
$ws.Range("M99").Value = -3.111100000000079
$ws.Range("N99").Value = -5199.6667
$ws.Range("H128").Value = 3533.3333
$ws.Range("I128").Value = 3533.3333
$ws.Range("K128").Value = 10599.9999
$ws.Range("M128").Value = -8109.999899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1282116.1
$ws.Range("I58").Value = 1716398.2
$ws.Range("J58").Value = 3396.3333
$ws.Range("K58").Value = 1716398.2
$ws.Range("L58").Value = 3396.3333
$ws.Range("M58").Value = -1716195.2
$ws.Range("N58").Value = -3802.3333
$ws.Range("H94").Value = 1233.25
$ws.Range("I94").Value = 912
$ws.Range("J94").Value = 1340.3334
$ws.Range("K94").Value = 912
$ws.Range("L94").Value = 1340.3334
$ws.Range("M94").Value = -461
$ws.Range("N94").Value = -2242.3334
$ws.Range("H105").Value = 1477.1428
$ws.Range("I105").Value = 868
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 868
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 879
$ws.Range("N105").Value = -6494
$ws.Range("H132").Value = 3772.9492
$ws.Range("I132").Value = 4044.4
$ws.Range("K132").Value = 12133.2
$ws.Range("M132").Value = -9603.200000000001
$ws.Range("H134").Value = 2860.6155
$ws.Range("I134").Value = 1829.5186
$ws.Range("J134").Value = 3974.2
$ws.Range("K134").Value = 5488.5558
$ws.Range("L134").Value = 11922.6
$ws.Range("M134").Value = -2953.5558
$ws.Range("N134").Value = -16992.6
$ws.Range("H136").Value = 1282116.1
$ws.Range("I136").Value = 1716398.2
$ws.Range("J136").Value = 3396.3333
$ws.Range("K136").Value = 5149194.6
$ws.Range("L136").Value = 10188.9999
$ws.Range("M136").Value = -5146644.6
$ws.Range("N136").Value = -15288.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 70.5
$ws.Range("I6").Value = 70.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 211.5
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -98.5
$ws.Range("H7").Value = 416.2
$ws.Range("I7").Value = 378.66666
$ws.Range("K7").Value = 1135.99998
$ws.Range("M7").Value = -1023.99998
$ws.Range("H11").Value = 3024.6667
$ws.Range("I11").Value = 298.16666
$ws.Range("J11").Value = 5751.1665
$ws.Range("K11").Value = 894.4999799999999
$ws.Range("L11").Value = 17253.4995
$ws.Range("M11").Value = -754.4999799999999
$ws.Range("N11").Value = -17533.4995
$ws.Range("H38").Value = 62.058823
$ws.Range("I38").Value = 35
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 105
$ws.Range("L38").Value = 450
$ws.Range("M38").Value = 242
$ws.Range("N38").Value = -1144
$ws.Range("H131").Value = 1484.025
$ws.Range("I131").Value = 2151.25
$ws.Range("J131").Value = 1317.2188
$ws.Range("K131").Value = 6453.75
$ws.Range("L131").Value = 3951.6564
$ws.Range("M131").Value = -1413.75
$ws.Range("N131").Value = -14031.6564

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2795.2083
$ws.Range("I126").Value = 1955.5555
$ws.Range("J126").Value = 3299
$ws.Range("K126").Value = 5866.666499999999
$ws.Range("L126").Value = 9897
$ws.Range("M126").Value = -3396.666499999999
$ws.Range("N126").Value = -14837
$ws.Range("H132").Value = 6228.773
$ws.Range("I132").Value = 4314.205
$ws.Range("J132").Value = 21162.4
$ws.Range("K132").Value = 12942.615
$ws.Range("L132").Value = 63487.2
$ws.Range("M132").Value = -10412.615
$ws.Range("N132").Value = -68547.20000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1845
$ws.Range("I22").Value = 1900
$ws.Range("J22").Value = 1826.6666
$ws.Range("K22").Value = 1900
$ws.Range("L22").Value = 1826.6666
$ws.Range("M22").Value = -1605
$ws.Range("N22").Value = -2416.6666
$ws.Range("H27").Value = 1845
$ws.Range("I27").Value = 1900
$ws.Range("J27").Value = 1826.6666
$ws.Range("K27").Value = 1900
$ws.Range("L27").Value = 1826.6666
$ws.Range("M27").Value = -1793
$ws.Range("N27").Value = -2040.6666
$ws.Range("H136").Value = 3062.3206
$ws.Range("I136").Value = 1751.2909
$ws.Range("J136").Value = 6197.391
$ws.Range("K136").Value = 5253.8727
$ws.Range("L136").Value = 18592.173
$ws.Range("M136").Value = -2703.8727
$ws.Range("N136").Value = -23692.173

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4190.6
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 7681.2
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 23043.6
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -26883.6
